$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.691.96'
$ws.Range('E2').Value = '  +0.85%  '

$ws.Range('D3').Value = '3.110.11'
$ws.Range('E3').Value = '  -0.53%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '585.14'
$ws.Range('E5').Value = '  -0.54%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.44'
$ws.Range('E6').Value = '  +0.45%  '

$ws.Range('E7').Value = '  +0.05%  '

$ws.Range('D8').Value = '3.105.03'
$ws.Range('E8').Value = '  -0.44%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.529'
$ws.Range('E9').Value = '  -0.28%  '

$ws.Range('E10').Value = '  +5.59%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.70'
$ws.Range('E11').Value = '  -0.50%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.460'
$ws.Range('E12').Value = '  -2.22%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000249'
$ws.Range('E13').Value = '  +0.10%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.41'
$ws.Range('E14').Value = '  +4.68%  '

$ws.Range('E15').Value = '  -0.54%  '

$ws.Range('D16').Value = '3.625.05'
$ws.Range('E16').Value = '  -0.49%  '

$ws.Range('D17').Value = '63.594.98'
$ws.Range('E17').Value = '  +0.85%  '

$ws.Range('E18').Value = '  -1.15%  '

$ws.Range('D19').Value = '3.102.08'
$ws.Range('E19').Value = '  -0.76%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '462.14'
$ws.Range('E20').Value = '  -0.82%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.33'
$ws.Range('E21').Value = '  +1.60%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.728'
$ws.Range('E22').Value = '  -0.72%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.55'
$ws.Range('E23').Value = '  -0.89%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.11'
$ws.Range('E24').Value = '  -3.25%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '81.43'
$ws.Range('E25').Value = '  -1.22%  '

$ws.Range('E26').Value = '  +0.12%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.00'
$ws.Range('E27').Value = '  +8.03%  '

$ws.Range('B28').Value = 'ImmutableX'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.22'
$ws.Range('E28').Value = '  -1.81%  '

$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.68'
$ws.Range('E29').Value = '  -1.52%  '

$ws.Range('E30').Value = '  +0.04%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.81'
$ws.Range('E31').Value = '  -0.90%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '26.88'
$ws.Range('E32').Value = '  -1.19%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.108'
$ws.Range('E33').Value = '  -3.08%  '

$ws.Range('D34').Value = '0.0₃0857'
$ws.Range('E34').Value = '  +5.04%  '

$ws.Range('B35').Value = 'Mantle'
$ws.Range('C35').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.04'
$ws.Range('E35').Value = '  +0.48%  '

$ws.Range('B36').Value = 'Stacks'
$ws.Range('C36').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.34'
$ws.Range('E36').Value = '  -2.53%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.40'
$ws.Range('E37').Value = '  +8.51%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.05'
$ws.Range('E38').Value = '  -0.37%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '50.85'
$ws.Range('E39').Value = '  -2.04%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '440.91'
$ws.Range('E40').Value = '  +2.63%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.70'
$ws.Range('E41').Value = '  -1.33%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0370'
$ws.Range('E42').Value = '  -0.95%  '

$ws.Range('D43').Value = '2.872.30'
$ws.Range('E43').Value = '  -2.79%  '

$ws.Range('E44').Value = '  -1.16%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.110'
$ws.Range('E45').Value = '  -2.08%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.15'
$ws.Range('E46').Value = '  -1.48%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '35.75'
$ws.Range('E47').Value = '  +1.97%  '

$ws.Range('E48').Value = '  +0.02%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '123.26'
$ws.Range('E49').Value = '  -1.71%  '

$ws.Range('E50').Value = '  -0.94%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '24.51'
$ws.Range('E51').Value = '  -2.23%  '
